$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D values are numeric-looking text (prices with "." as thousands separators
# or with meaningful trailing zeros). Force text format so Excel does not coerce them
# into floating point numbers and lose formatting/precision.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '57.055.33'
$ws.Range("E2").Value = '  +1.11%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.434.11'
$ws.Range("E3").Value = '  -1.71%  '
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '489.94'
$ws.Range("E5").Value = '  +0.29%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '154.75'
$ws.Range("E6").Value = '  +4.12%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.618'
$ws.Range("E7").Value = '  +21.26%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.996'
$ws.Range("E8").Value = '  -0.18%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.450.72'
$ws.Range("E9").Value = '  -1.43%  '
$ws.Range("E10").Value = '  +2.51%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.63'
$ws.Range("E11").Value = '  -2.58%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.338'
$ws.Range("E12").Value = '  +1.40%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.126'
$ws.Range("E13").Value = '  +1.00%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '2.859.17'
$ws.Range("E14").Value = '  -1.76%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '57.147.35'
$ws.Range("E15").Value = '  +1.32%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '20.80'
$ws.Range("E16").Value = '  -0.76%  '
$ws.Range("E17").Value = '  -1.12%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.445.95'
$ws.Range("E18").Value = '  -2.10%  '
$ws.Range("E19").Value = '  +5.73%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '329.22'
$ws.Range("E20").Value = '  +3.59%  '
$ws.Range("E21").Value = '  -1.50%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.00'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.93'
$ws.Range("E23").Value = '  +1.10%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '58.08'
$ws.Range("E24").Value = '  -0.64%  '
$ws.Range("E25").Value = '  +0.61%  '
$ws.Range("E26").Value = '  -0.15%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.160'
$ws.Range("E27").Value = '  -1.56%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.548.40'
$ws.Range("E28").Value = '  -1.92%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.39'
$ws.Range("E29").Value = '  -2.60%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0₃0791'
$ws.Range("E30").Value = '  +0.63%  '
$ws.Range("E31").Value = '  -0.04%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '18.80'
$ws.Range("E32").Value = '  +3.02%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '149.78'
$ws.Range("E33").Value = '  +0.62%  '
$ws.Range("E34").Value = '  +1.27%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.32'
$ws.Range("E35").Value = '  +2.66%  '
$ws.Range("E36").Value = '  -0.76%  '
$ws.Range("E37").Value = '  -0.78%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.856'
$ws.Range("E38").Value = '  -1.48%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.103'
$ws.Range("E39").Value = '  +11.45%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '34.19'
$ws.Range("E40").Value = '  +1.24%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.38'
$ws.Range("E41").Value = '  +0.42%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.53'
$ws.Range("E42").Value = '  +1.08%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.994'
$ws.Range("E43").Value = '  -0.18%  '
$ws.Range("E44").Value = '  -1.62%  '
$ws.Range("E45").Value = '  -2.95%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '267.68'
$ws.Range("E46").Value = '  +1.42%  '
$ws.Range("B47").Value = 'VeChain'
$ws.Range("C47").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0230'
$ws.Range("E47").Value = '  +0.21%  '
$ws.Range("B48").Value = 'WhiteBITCoin'
$ws.Range("C48").Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '10.20'
$ws.Range("E48").Value = '  -0.25%  '
$ws.Range("B49").Value = 'RenderToken'
$ws.Range("C49").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '4.68'
$ws.Range("E49").Value = '  -1.94%  '
$ws.Range("E50").Value = '  +0.10%  '
$ws.Range("B51").Value = 'Maker'
$ws.Range("C51").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.862.95'
$ws.Range("E51").Value = '  -1.40%  '
